$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "id_hospital"
$null = $ws.Range("D13").Select()
